$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert a new row after the "Create a viewing area..." row (row 8), shifting
# rows 9-12 down to 10-13.
$ws.Rows.Item(9).Insert()

# New row 9: task for viewing CE info from MainWin (JDialog), still In Progress.
# Write the D/F text first, then the status cells, so the new shared strings
# land in the same order they do in the target workbook.
$ws.Range("D9").Value = "Create a JDialog to view CE information from MainWin"
$ws.Range("E8").Value = "Complete "
$ws.Range("F9").Value = "Currently trying to format the Jdialog"
$ws.Range("E9").Value = "In Progress"
$ws.Range("E10").Value = "Not Started"
$ws.Range("B9").Value = "LCE"
$ws.Range("E11").Value = "Not Started"

# Reflect the active selection left behind by the editor.
$ws.Range("D17").Select()
